$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs",  "Cxcl13", "Ackr4", "FAPs", 1, 0.3333333333333333, 0.2087556666666667, 0.626267, 0.01876624903294638, 0.01876624903294638, 3, 1, 1.890486333333333, 5.671459, 0.9442400689667343, 0.9442400689667344, 0.3946497348392223, 3.551847613553, 0.0177198442811162, 0.0177198442811162),
    @("ECs",  "Cxcl13", "Ackr4", "sCs",  1, 0.3333333333333333, 0.2087556666666667, 0.626267, 0.01876624903294638, 0.01876624903294638, 2, 0.6666666666666666, 0.1116383333333333, 0.334915, 0.05575993103326565, 0.05575993103326566, 0.02330513470055556, 0.209746212305, 0.001046404751830178, 0.001046404751830178),
    @("FAPs", "Cxcl13", "Ackr4", "FAPs", 3, 1, 10.91523966666667, 32.745719, 0.9812337509670536, 0.9812337509670537, 3, 1, 1.890486333333333, 5.671459, 0.9442400689667343, 0.9442400689667344, 20.63511141489122, 185.716002734021, 0.926520224685618, 0.9265202246856182),
    @("FAPs", "Cxcl13", "Ackr4", "sCs",  3, 1, 10.91523966666667, 32.745719, 0.9812337509670536, 0.9812337509670537, 2, 0.6666666666666666, 0.1116383333333333, 0.334915, 0.05575993103326565, 0.05575993103326566, 1.218559164320556, 10.967032478885, 0.05471352628143547, 0.05471352628143548)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $colNum = $c + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$c]
    }
}
